$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(4)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 3 of the body placeholder: “Hello, world.” for a string literal
$para = $tr.Paragraphs(3)

# --- Step 1: opening-quote run -> append a straight quote after the curly quote ---
# (position 1, 1 char: the opening curly quote)
$openQuote = $para.Characters(1, 1)
$openQuote.Text = [char]0x201C + '"'

# Re-fetch paragraph (text grew by 1 char after step 1).
$para = $tr.Paragraphs(3)

# --- Step 2: split "Hello, world." into "Hello, " and "world." ---
# "Hello, " is now at positions 3-9 (after the 2-char opening-quote run), "world." at 10-15.
$worldPart = $para.Characters(10, 6)
$worldPart.Text = "world."

# --- Step 3: insert a straight quote before the closing curly quote, keeping ---
# --- it and the following space together as one run: `"” ` ---
$para = $tr.Paragraphs(3)
# the closing curly quote + trailing space are now at positions 16-17
$closeQuotePlusSpace = $para.Characters(16, 2)
$closeQuotePlusSpace.Text = '"' + [char]0x201D + ' '

# --- Step 4: split off the trailing "for a string literal" text ---
$para = $tr.Paragraphs(3)
# after step 3 the paragraph reads: “"Hello, world."”  for a string literal
$tail = $para.Characters(19, 20)
$tail.Text = "for a string literal"
